$wb = $excel.ActiveWorkbook

# --- Sheet "login Details" ---
$ws1 = $wb.Worksheets.Item("login Details")

# A2: standard_user -> standard_userNkosi
$ws1.Range("A2").Value = "standard_userNkosi"
$ws1.Columns.Item(1).AutoFit()

# Clear column B contents that held duplicate "secret_sauce" / "Test" entries (rows 3-7, 9-10)
$ws1.Range("B3").ClearContents()
$ws1.Range("B4").ClearContents()
$ws1.Range("B5").ClearContents()
$ws1.Range("B6").ClearContents()
$ws1.Range("B7").ClearContents()
$ws1.Range("B9").ClearContents()
$ws1.Range("B10").ClearContents()

# Update selection to F7
$ws1.Range("F7").Select()

$wb.Save()
